$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.836.49'
$ws.Range("E2").Value = '  -2.28%  '
$ws.Range("D3").Value = '''3.543.38'
$ws.Range("E3").Value = '  -3.40%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''616.44'
$ws.Range("E5").Value = '  -4.04%  '
$ws.Range("D6").Value = '''154.09'
$ws.Range("E6").Value = '  -2.84%  '
$ws.Range("D7").Value = '''3.538.82'
$ws.Range("E7").Value = '  -3.49%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '''0.487'
$ws.Range("E9").Value = '  -1.74%  '
$ws.Range("E10").Value = '  -2.08%  '
$ws.Range("D11").Value = '''6.94'
$ws.Range("E11").Value = '  -1.41%  '
$ws.Range("D12").Value = '''0.432'
$ws.Range("E12").Value = '  -1.45%  '
$ws.Range("E13").Value = '  -2.05%  '
$ws.Range("D14").Value = '''32.25'
$ws.Range("E14").Value = '  -0.03%  '
$ws.Range("D15").Value = '''4.141.42'
$ws.Range("E15").Value = '  -3.51%  '
$ws.Range("D16").Value = '''3.539.02'
$ws.Range("E16").Value = '  -3.49%  '
$ws.Range("D17").Value = '''67.768.52'
$ws.Range("E17").Value = '  -2.45%  '
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("D19").Value = '''6.39'
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("D21").Value = '''453.99'
$ws.Range("E21").Value = '  -2.47%  '
$ws.Range("E22").Value = '  -2.70%  '
$ws.Range("D23").Value = '''0.643'
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").Value = '''78.09'
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("D25").Value = '''3.683.82'
$ws.Range("E25").Value = '  -3.44%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '''0.0000120'
$ws.Range("E27").Value = '  -3.40%  '
$ws.Range("D28").Value = '''10.52'
$ws.Range("E28").Value = '  -2.29%  '
$ws.Range("E29").Value = '  -5.33%  '
$ws.Range("D30").Value = '''2.57'
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("D31").Value = '''1.66'
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("E33").Value = '  -2.08%  '
$ws.Range("E34").Value = '  -3.62%  '
$ws.Range("D35").Value = '''6.24'
$ws.Range("E35").Value = '  -2.62%  '
$ws.Range("E36").Value = '  -2.19%  '
$ws.Range("D37").Value = '''3.540.75'
$ws.Range("E37").Value = '  -3.29%  '
$ws.Range("E38").Value = '  -3.53%  '
$ws.Range("D40").Value = '''0.999'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").Value = '''176.12'
$ws.Range("E41").Value = '  -1.87%  '
$ws.Range("E42").Value = '  -4.20%  '
$ws.Range("D43").Value = '''0.0882'
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("D44").Value = '''2.10'
$ws.Range("E44").Value = '  -3.42%  '
$ws.Range("E45").Value = '  -4.21%  '
$ws.Range("D46").Value = '''29.55'
$ws.Range("E46").Value = '  +10.10%  '
$ws.Range("E47").Value = '  -1.70%  '
$ws.Range("D48").Value = '''2.60'
$ws.Range("E48").Value = '  -3.16%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Value = '''1.23'
$ws.Range("E49").Value = '  -1.61%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '''7.69'
$ws.Range("E50").Value = '  -1.18%  '
$ws.Range("E51").Value = '  -3.00%  '
